# Sales report "models clean up": refresh order rows with the 02/24/2023
# batch, renumbering order IDs and recording the new totals/payment info,
# extending the sheet down to row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds order-date text like "02/12/2023" stored as literal text,
# not a real date. Force text formatting first so Excel doesn't reinterpret
# the "02/24/2023" strings we assign as date serials.
$ws.Range("A2:A16").NumberFormat = "@"

$rows = @(
    @{ Row=2;  Date="02/24/2023"; OrderId=1;  Product="";        Price=0;     Total=180000; Payment="COD";       Status="pending" },
    @{ Row=3;  Date="02/24/2023"; OrderId=2;  Product="";        Price=0;     Total=120000; Payment="Razor Pay"; Status="pending" },
    @{ Row=4;  Date="02/24/2023"; OrderId=3;  Product="";        Price=0;     Total=120000; Payment="Wallet";    Status="pending" },
    @{ Row=5;  Date="02/24/2023"; OrderId=4;  Product="";        Price=0;     Total=300000; Payment="COD";       Status="pending" },
    @{ Row=6;  Date="02/24/2023"; OrderId=5;  Product="Oppo x3"; Price=60000; Total=60000;  Payment="Wallet";    Status="pending" },
    @{ Row=7;  Date="02/24/2023"; OrderId=6;  Product="";        Price=0;     Total=60000;  Payment="COD";       Status="pending" },
    @{ Row=8;  Date="02/24/2023"; OrderId=7;  Product="";        Price=0;     Total=60000;  Payment="COD";       Status="pending" },
    @{ Row=9;  Date="02/24/2023"; OrderId=8;  Product="";        Price=0;     Total=120000; Payment="COD";       Status="pending" },
    @{ Row=10; Date="02/24/2023"; OrderId=9;  Product="";        Price=0;     Total=60000;  Payment="COD";       Status="pending" },
    @{ Row=11; Date="02/24/2023"; OrderId=10; Product="";        Price=0;     Total=60000;  Payment="COD";       Status="pending" },
    @{ Row=12; Date="02/24/2023"; OrderId=11; Product="";        Price=0;     Total=60000;  Payment="COD";       Status="pending" },
    @{ Row=13; Date="02/24/2023"; OrderId=13; Product="";        Price=0;     Total=60000;  Payment="Wallet";    Status="pending" },
    @{ Row=14; Date="02/24/2023"; OrderId=12; Product="";        Price=0;     Total=60000;  Payment="Wallet";    Status="pending" },
    @{ Row=15; Date="02/24/2023"; OrderId=14; Product="";        Price=0;     Total=0;      Payment="";          Status="" },
    @{ Row=16; Date="02/24/2023"; OrderId=15; Product="";        Price=0;     Total=0;      Payment="";          Status="" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.OrderId
    $ws.Cells.Item($row, 3).Value = $r.Product
    $ws.Cells.Item($row, 4).Value = $r.Price
    $ws.Cells.Item($row, 5).Value = $r.Total
    $ws.Cells.Item($row, 6).Value = $r.Payment
    $ws.Cells.Item($row, 7).Value = $r.Status
}
